$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-9, columns A-T
# Column order: A B C D E F G H I J K L M N O P Q R S T
$data = @(
    @{ row = 2;  A = "MuSCs";         B = "Ifng"; C = "Ifngr1"; D = "ECs";           E = 1; F = 0.3333333333333333; G = 0.06628299999999999;  H = 0.198849; I = 0.581839834503261; J = 0.581839834503261; K = 3; L = 1; M = 76.56242133333333; N = 229.687264; O = 0.3581170284509773;  P = 0.3581170284509773;  Q = 5.074786973237333;  R = 45.673082759136;      S = 0.2083667525667163;  T = 0.2083667525667162 },
    @{ row = 3;  A = "MuSCs";         B = "Ifng"; C = "Ifngr1"; D = "FAPs";          E = 1; F = 0.3333333333333333; G = 0.06628299999999999;  H = 0.198849; I = 0.581839834503261; J = 0.581839834503261; K = 3; L = 1; M = 47.437018;         N = 142.311054; O = 0.221884360876912;   P = 0.221884360876912;   Q = 3.144267864094;     R = 28.298410776846;      S = 0.1291011598114843;  T = 0.1291011598114843 },
    @{ row = 4;  A = "MuSCs";         B = "Ifng"; C = "Ifngr1"; D = "MuSCs";         E = 1; F = 0.3333333333333333; G = 0.06628299999999999;  H = 0.198849; I = 0.581839834503261; J = 0.581839834503261; K = 3; L = 1; M = 18.65483866666667; N = 55.964516;  O = 0.08725710698794852; P = 0.08725710698794852; Q = 1.236498671342667;  R = 11.128488042084;      S = 0.0507696606891013;  T = 0.0507696606891013 },
    @{ row = 5;  A = "MuSCs";         B = "Ifng"; C = "Ifngr1"; D = "Resolving-Mac"; E = 1; F = 0.3333333333333333; G = 0.06628299999999999;  H = 0.198849; I = 0.581839834503261; J = 0.581839834503261; K = 3; L = 1; M = 71.13734666666666; N = 213.41204;  O = 0.3327415036841621;  P = 0.3327415036841621;  Q = 4.715196749106666;  R = 42.43677074196;       S = 0.1936022614359591;  T = 0.1936022614359591 },
    @{ row = 6;  A = "Resolving-Mac"; B = "Ifng"; C = "Ifngr1"; D = "ECs";           E = 1; F = 0.3333333333333333; G = 0.04763666666666667; H = 0.14291;  I = 0.4181601654967389; J = 0.4181601654967389; K = 3; L = 1; M = 76.56242133333333; N = 229.687264; O = 0.3581170284509773;  P = 0.3581170284509773;  Q = 3.647178544248889;  R = 32.82460689824;       S = 0.149750275884261;   T = 0.149750275884261 },
    @{ row = 7;  A = "Resolving-Mac"; B = "Ifng"; C = "Ifngr1"; D = "FAPs";          E = 1; F = 0.3333333333333333; G = 0.04763666666666667; H = 0.14291;  I = 0.4181601654967389; J = 0.4181601654967389; K = 3; L = 1; M = 47.437018;         N = 142.311054; O = 0.221884360876912;   P = 0.221884360876912;   Q = 2.259741414126667;  R = 20.33767272714;       S = 0.09278320106542767; T = 0.09278320106542767 },
    @{ row = 8;  A = "Resolving-Mac"; B = "Ifng"; C = "Ifngr1"; D = "MuSCs";         E = 1; F = 0.3333333333333333; G = 0.04763666666666667; H = 0.14291;  I = 0.4181601654967389; J = 0.4181601654967389; K = 3; L = 1; M = 18.65483866666667; N = 55.964516;  O = 0.08725710698794852; P = 0.08725710698794852; Q = 0.8886543312844445; R = 7.997888981560001;    S = 0.0364874462988472;  T = 0.0364874462988472 },
    @{ row = 9;  A = "Resolving-Mac"; B = "Ifng"; C = "Ifngr1"; D = "Resolving-Mac"; E = 1; F = 0.3333333333333333; G = 0.04763666666666667; H = 0.14291;  I = 0.4181601654967389; J = 0.4181601654967389; K = 3; L = 1; M = 71.13734666666666; N = 213.41204;  O = 0.3327415036841621;  P = 0.3327415036841621;  Q = 3.388746070711111;  R = 30.4987146364;        S = 0.139139242248203;   T = 0.139139242248203 }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($entry in $data) {
    $r = $entry.row
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $entry[$col]
    }
}
